# Homework 2.pptx edit — "sent data to sql"
#
# 1) The cached "datetimeFigureOut" footer field (2021/10/27 -> 2021/10/28)
#    on the slide master and on every slide layout.
# 2) Strike-through the "報名資訊存入資料庫 +10" bullet on slide 7 (the
#    bonus-points slide), marking that item as already done ("sent data
#    to sql").

$p = $ppt.ActivePresentation

$oldDate = "2021/10/27"
$newDate = "2021/10/28"

# --- 1) Update the date placeholder everywhere it is cached -------------

$master = $p.SlideMaster

# Slide master itself.
for ($k = 1; $k -le $master.Shapes.Count; $k++) {
    $msh = $master.Shapes.Item($k)
    if ($msh.HasTextFrame -and $msh.TextFrame.HasText -and $msh.TextFrame.TextRange.Text -eq $oldDate) {
        $msh.TextFrame.TextRange.Text = $newDate
    }
}

# Every slide layout hanging off the master.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lay = $layouts.Item($i)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $sh = $lay.Shapes.Item($j)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2) Strike through the "報名資訊存入資料庫 +10" bullet on slide 7 ----

$target = "報名資訊存入資料庫"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($sj = 1; $sj -le $slide.Shapes.Count; $sj++) {
        $shape = $slide.Shapes.Item($sj)
        if (-not $shape.HasTextFrame) { continue }
        if (-not $shape.TextFrame.HasText) { continue }

        $tr = $shape.TextFrame.TextRange
        $count = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $count; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            if ($para.Text.Contains($target)) {
                $para.Font.Strike = $true
            }
        }
    }
}

Write-Output "done"
